$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197; existing rows 197-273 shift down to 198-274.
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row 197 with the new record
$ws.Range("A197").Value = 5
$ws.Range("B197").Value = "Macroferia Regional de Talca"
$ws.Range("C197").Value = "Maule"
$ws.Range("D197").Value = 45027
$ws.Range("E197").Value = 7
$ws.Range("F197").Value = 100112017
$ws.Range("G197").Value = "Apio"
$ws.Range("H197").Value = "Americana (o)"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 600
$ws.Range("K197").Value = 8000
$ws.Range("L197").Value = 8000
$ws.Range("M197").Value = 8000
$ws.Range("N197").Value = "`$/docena de matas"
$ws.Range("O197").Value = "Provincia del Elquí"
$ws.Range("P197").Value = 1333
$ws.Range("Q197").Value = 6
$ws.Range("R197").Value = "Hortaliza"
